$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last data row (row 10) - table shrinks from 10 data rows to 9
$ws.Rows.Item(10).Delete()

# New table data (ESPN comparison for NBA, Wednesday 6th Mar 2024)
$data = @(
    @("NBA, Wednesday 6th Mar 2024", "Ballgorithm", "ESPN"),
    @("Orlando Magic (36-26) vs Washington Wizards (9-52)", "Orlando Magic (72.41%)", "Orlando Magic (69.7%)"),
    @("Cleveland Cavaliers (40-21) vs Atlanta Hawks (27-34)", "Cleveland Cavaliers (65.62%)", "Cleveland Cavaliers (62.4%)"),
    @("Los Angeles Clippers (39-21) vs Houston Rockets (27-34)", "Los Angeles Clippers (72.41%)", "Los Angeles Clippers (70.2%)"),
    @("Memphis Grizzlies (21-41) vs Philadelphia 76ers (35-26)", "Philadelphia 76ers (59.38%)", "Philadelphia 76ers (79.2%)"),
    @("Chicago Bulls (29-32) vs Utah Jazz (28-34)", "Utah Jazz (63.33%)", "Chicago Bulls (50.5%)"),
    @("Milwaukee Bucks (41-21) vs Golden State Warriors (32-28)", "Milwaukee Bucks (78.12%)", "Golden State Warriors (59.1%)"),
    @("Oklahoma City Thunder (42-19) vs Portland Trail Blazers (17-43)", "Oklahoma City Thunder (80.00%)", "Oklahoma City Thunder (87.4%)"),
    @("Sacramento Kings (34-26) vs Los Angeles Lakers (35-29)", "Los Angeles Lakers (69.70%)", "Los Angeles Lakers (54.9%)")
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}

$ws.Range("C10").Select()
